$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Shallow water port")

# Reassign ship names in column B (rows 4-21) to reflect the corrected
# shared-string order, and correct the BR values in column C so that
# each ship keeps its proper BR rating.
$names = @("Hercules", "Pandora", "Mercury", "Mortar Brig", "NavyBrig", "Niagara", "Prince de Neufchatel", "Rattlesnake", "Rattlesnake Heavy", "Snow", "Brig", "Pickle", "Cutter", "GunBoat", "Lynx", "Privateer", "Yacht", "Yacht Silver")
$brs = @(100, 100, 80, 80, 80, 80, 80, 80, 80, 80, 70, 55, 50, 50, 50, 50, 50, 50)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 2).Value = $names[$i]
    $ws.Cells.Item($row, 3).Value = $brs[$i]
}

# Extend the BR/player sums to include the last row (21) which was
# previously excluded from the totals.
$ws.Range("D3").Formula = "=SUM(D4:D21)"
$ws.Range("E3").Formula = "=SUM(E4:E21)"

$ws1 = $wb.Worksheets.Item("Deep water port")
$ws1.Range("D3").Formula = "=SUM(D4:D35)"
$ws1.Range("E3").Formula = "=SUM(E4:E35)"
